# Updated cryptos list - applies price/volume changes from the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.463.74'
$ws.Range("E2").Value = '  +0.79%  '

# Row 3
$ws.Range("D3").Value = '3.531.43'
$ws.Range("E3").Value = '  +0.55%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").Value = '''597.46'
$ws.Range("E5").Value = '  +0.82%  '

# Row 6
$ws.Range("D6").Value = '''173.75'
$ws.Range("E6").Value = '  +2.11%  '

# Row 7
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("E8").Value = '  +2.62%  '

# Row 9
$ws.Range("E9").Value = '  +8.36%  '

# Row 10
$ws.Range("D10").Value = '''7.31'
$ws.Range("E10").Value = '  +0.78%  '

# Row 11
$ws.Range("E11").Value = '  -0.20%  '

# Row 12
$ws.Range("D12").Value = '4.143.82'
$ws.Range("E12").Value = '  +0.83%  '

# Row 13
$ws.Range("E13").Value = '  -0.12%  '

# Row 14
$ws.Range("D14").Value = '''28.89'
$ws.Range("E14").Value = '  +2.80%  '

# Row 15
$ws.Range("E15").Value = '  +2.61%  '

# Row 16
$ws.Range("D16").Value = '67.408.65'
$ws.Range("E16").Value = '  +0.83%  '

# Row 17
$ws.Range("D17").Value = '3.533.77'
$ws.Range("E17").Value = '  +1.24%  '

# Row 18
$ws.Range("E18").Value = '  +1.23%  '

# Row 19
$ws.Range("D19").Value = '''14.23'
$ws.Range("E19").Value = '  +1.57%  '

# Row 20
$ws.Range("D20").Value = '''398.89'
$ws.Range("E20").Value = '  +2.29%  '

# Row 21
$ws.Range("E21").Value = '  +0.67%  '

# Row 22
$ws.Range("D22").Value = '''73.60'
$ws.Range("E22").Value = '  +0.52%  '

# Row 23
$ws.Range("E23").Value = '  +2.43%  '

# Row 24
$ws.Range("D24").Value = '''0.998'
$ws.Range("E24").Value = '  -0.26%  '

# Row 25
$ws.Range("E25").Value = '  -0.25%  '

# Row 26
$ws.Range("E26").Value = '  +1.65%  '

# Row 27
$ws.Range("E27").Value = '  +0.29%  '

# Row 28
$ws.Range("E28").Value = '  -0.09%  '

# Row 29
$ws.Range("D29").Value = '''6.30'
$ws.Range("E29").Value = '  -1.36%  '

# Row 30
$ws.Range("D30").Value = '''1.47'
$ws.Range("E30").Value = '  -0.01%  '

# Row 31
$ws.Range("E31").Value = '  +1.25%  '

# Row 32
$ws.Range("E32").Value = '  +2.65%  '

# Row 33
$ws.Range("E33").Value = '  +0.10%  '

# Row 34
$ws.Range("E34").Value = '  +3.89%  '

# Row 35
$ws.Range("D35").Value = '''164.22'
$ws.Range("E35").Value = '  +1.87%  '

# Row 36
$ws.Range("D36").Value = '''0.897'
$ws.Range("E36").Value = '  -0.62%  '

# Row 37
$ws.Range("E37").Value = '  -0.74%  '

# Row 38
$ws.Range("D38").Value = '''6.99'
$ws.Range("E38").Value = '  +3.93%  '

# Row 39
$ws.Range("D39").Value = '''4.76'
$ws.Range("E39").Value = '  +1.90%  '

# Row 40
$ws.Range("D40").Value = '''0.0750'
$ws.Range("E40").Value = '  +0.28%  '

# Row 41
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").Value = '''27.40'
$ws.Range("E41").Value = '  +2.22%  '

# Row 42
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = '''26.62'
$ws.Range("E42").Value = '  +0.56%  '

# Row 43
$ws.Range("D43").Value = '''2.64'
$ws.Range("E43").Value = '  +3.22%  '

# Row 44
$ws.Range("D44").Value = '2.809.61'
$ws.Range("E44").Value = '  +0.17%  '

# Row 45
$ws.Range("D45").Value = '''42.98'
$ws.Range("E45").Value = '  -1.24%  '

# Row 46
$ws.Range("D46").Value = '''0.0312'
$ws.Range("E46").Value = '  -0.62%  '

# Row 47
$ws.Range("D47").Value = '''342.38'
$ws.Range("E47").Value = '  -3.83%  '

# Row 48
$ws.Range("E48").Value = '  +0.97%  '

# Row 49
$ws.Range("D49").Value = '''33.98'
$ws.Range("E49").Value = '  +2.16%  '

# Row 50
$ws.Range("D50").Value = '''6.55'
$ws.Range("E50").Value = '  +0.88%  '

# Row 51
$ws.Range("E51").Value = '  +0.59%  '
